# Apply the "Increased elysis in FR00 2040" edit to the Capacity sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Capacity")

# --- Update capacity values (Elec_capa = col F, Other_capa = col H) ---
$ws.Range("F42").Value = 140   # FI00 Hydrogen storage dimensioner 2040
$ws.Range("F47").Value = 120   # SE02 Hydrogen storage dimensioner 2040
$ws.Range("H48").Value = 2000  # SE02 Hydrogen processor 2040
$ws.Range("F50").Value = 80    # NON1 Hydrogen storage dimensioner 2040
$ws.Range("F66").Value = 300   # DE00 Hydrogen storage dimensioner 2040
$ws.Range("H80").Value = 15000 # FR00 Electrolysis 2040

# --- Reveal "Hydrogen storage dimensioner" rows for 2040 by widening the
#     AutoFilter on the Generator_ID column (col C, 0-based colId 2 / field 3)
#     to include that value alongside the existing ones. ---
$ws.Range("A1:J177").AutoFilter(3, @("Electrolysis", "Hydrogen processor", "Hydrogen storage dimensioner"), 7)

# --- Update the active selection to reflect the last edited cell ---
$ws.Activate()
$ws.Range("H81").Select()
